$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '42.815.77'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  -1.36%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.334.41'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +0.87%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.09%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '306.72'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  -1.69%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '100.23'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -2.26%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.508'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  -5.33%  '
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.512'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  -3.78%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '34.90'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  -2.84%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '52.03'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  -0.69%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0798'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  -1.98%  '
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -0.80%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.80'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  -3.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '15.54'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.33%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.350.27'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +2.41%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.794'
$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  -2.20%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '42.718.70'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -1.37%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.28'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +1.90%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.0₃0904'
$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -1.93%  '
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '11.57'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  -7.32%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '67.57'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  -1.30%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '236.07'
$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -2.63%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.99'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  -2.74%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.56'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  -2.74%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.999'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.22%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '24.88'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.63%  '
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +3.22%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.82'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -7.85%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '9.34'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  -3.24%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '158.78'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  -4.90%  '
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  -0.15%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.10'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  -4.17%  '
$ws.Range('B34').Value = 'WEMIXToken'
$ws.Range('C34').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.45'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  -4.12%  '
$ws.Range('B35').Value = 'Celestia'
$ws.Range('C35').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.32'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  -2.65%  '
$ws.Range('B36').Value = 'Hedera'
$ws.Range('C36').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0724'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  -2.88%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.55'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +5.04%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.95'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  -4.68%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.86'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +0.03%  '
$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  -4.54%  '
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  -3.79%  '
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  +1.90%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.007.05'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +1.22%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.0284'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  -1.62%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '18.66'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -5.83%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '10.29'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +4.65%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.91'
$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.47%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '55.66'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -0.94%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.92'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  -1.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.556.36'
$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  +0.56%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '4.66'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  +3.07%  '
